# Generate Report for Handback
# Adds a new handback entry (4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md) as row 4
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet (sheet1) - new row 4
# ---------------------------------------------------------------------------
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value() = "4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aff22a4bb7b1a6361dbc618bcce03b82fb8d2aa3/e2e/4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md", "", "", "e2e\4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md") | Out-Null
$wsOverview.Range("C4").Value() = ".md"
$wsOverview.Range("E4").Value() = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value() = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value() = "2016-09-05 18:50:17"

# ---------------------------------------------------------------------------
# zh-cn sheet (sheet2) - new row 4
# ---------------------------------------------------------------------------
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aff22a4bb7b1a6361dbc618bcce03b82fb8d2aa3/e2e/4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md", "", "", "4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md") | Out-Null
$wsZhCn.Range("B4").Value() = ".md"
$wsZhCn.Range("C4").Value() = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value() = "e2e"
$wsZhCn.Range("E4").Value() = "ht"
$wsZhCn.Range("F4").Value() = "'True"
$wsZhCn.Range("G4").Value() = "4a94e237-c3b1-42cf-ad8e-648b4c7ab206.d805dbdef0e2659a71af577ae30c7f6123b29d33.zh-cn.xlf"
$wsZhCn.Range("H4").Value() = "2016-09-05 18:50:11"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/645a532a337a64ac81f98fbbf0c4f8daf2d5e3f5/e2e/4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md", "", "", "4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md") | Out-Null
$wsZhCn.Range("J4").Value() = "4a94e237-c3b1-42cf-ad8e-648b4c7ab206.d805dbdef0e2659a71af577ae30c7f6123b29d33.zh-cn.xlf"
$wsZhCn.Range("K4").Value() = "2016-09-05 18:50:35"
$wsZhCn.Range("L4").Value() = "'"
$wsZhCn.Range("M4").Value() = "'True"
$wsZhCn.Range("N4").Value() = "'"
$wsZhCn.Range("O4").Value() = "'False"
$wsZhCn.Range("P4").Value() = "'"

# ---------------------------------------------------------------------------
# de-de sheet (sheet3) - new row 4
# ---------------------------------------------------------------------------
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aff22a4bb7b1a6361dbc618bcce03b82fb8d2aa3/e2e/4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md", "", "", "4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md") | Out-Null
$wsDeDe.Range("B4").Value() = ".md"
$wsDeDe.Range("C4").Value() = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value() = "e2e"
$wsDeDe.Range("E4").Value() = "ht"
$wsDeDe.Range("F4").Value() = "'True"
$wsDeDe.Range("G4").Value() = "4a94e237-c3b1-42cf-ad8e-648b4c7ab206.d805dbdef0e2659a71af577ae30c7f6123b29d33.de-de.xlf"
$wsDeDe.Range("H4").Value() = "2016-09-05 18:50:17"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/26183a743785e74898db56092d288b2a468f9db0/e2e/4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md", "", "", "4a94e237-c3b1-42cf-ad8e-648b4c7ab206.md") | Out-Null
$wsDeDe.Range("J4").Value() = "4a94e237-c3b1-42cf-ad8e-648b4c7ab206.d805dbdef0e2659a71af577ae30c7f6123b29d33.de-de.xlf"
$wsDeDe.Range("K4").Value() = "2016-09-05 18:50:42"
$wsDeDe.Range("L4").Value() = "'"
$wsDeDe.Range("M4").Value() = "'True"
$wsDeDe.Range("N4").Value() = "'"
$wsDeDe.Range("O4").Value() = "'False"
$wsDeDe.Range("P4").Value() = "'"

Write-Host "Row 4 added to Overview, zh-cn and de-de sheets."
